$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.764084815979004
$ws.Range("B1").Value = 5.376907825469971
$ws.Range("C1").Value = 3.547969341278076
$ws.Range("D1").Value = 0.9269639253616333
$ws.Range("E1").Value = 0.5927579998970032
